# Applies the commit's changes:
#  1. Slide 2 ("The NYPD arrested 210,000 persons in 2019"): update the
#     "Our result" sentence with new precision/recall/F1 figures and remove
#     the following paragraph about predicting race from arrests.
#  2. Delete the slide "Brooklyn leads in the number of felonies " (slide 8
#     in the original running order).

$p = $ppt.ActivePresentation

# --- 1. Edit slide 2's body text -------------------------------------------------
$s2 = $p.Slides.Item(2)
$contentShape = $s2.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange

# Paragraph 4 is: "Our result" (styled run) + ": No, we could not predict ... (0.56 accuracy)"
# Only the second run's text is replaced, so the "Our result" run keeps its own formatting.
$para4 = $tr.Paragraphs(4, 1)
$firstRunLength = "Our result".Length
$secondRunStart = $para4.Start + $firstRunLength
$secondRunLength = $para4.Length - $firstRunLength
$secondRun = $tr.Characters($secondRunStart, $secondRunLength)
$secondRun.Text = ":  (0.56 precision) (0.85 recall) (0.69 F1 score) for Black"

# Paragraph 5 ("If these predictions cannot be improved ... do not predict on race.") is removed entirely.
$para5 = $tr.Paragraphs(5, 1)
$para5.Delete()

# --- 2. Delete the "Brooklyn leads in the number of felonies" slide --------------
$targetSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    foreach ($shp in $candidate.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "Brooklyn leads in the number of felonies ") {
                $targetSlide = $candidate
            }
            break
        }
    }
    if ($targetSlide -ne $null) { break }
}
if ($targetSlide -ne $null) {
    $targetSlide.Delete()
} else {
    # Fallback: in the original deck this slide is always at position 8.
    $p.Slides.Item(8).Delete()
}
